$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column E (shifts old E:N -> H:Q)
$ws.Range("E1:G1").EntireColumn.Insert()

# New Month/Day/Year columns get the same width as column D (9.5 chars)
$ws.Range("E1:G1").ColumnWidth = 8.67

# Fill in the new header cells (order matters for shared-string table ordering)
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"
$ws.Range("E1").Value = "Month"

# Fill in Month/Day/Year values for each data row (date sampled was 7/14/2015)
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = 2015

$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 2015

$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 2015

$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 2015

# Match the final selection from the authored edit
$ws.Range("E1:G1").Select() | Out-Null
